$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text formatting (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 17 and 18 swap coin name/link along with updated price/volume
$ws.Range("D2").Value = '27.331.36'
$ws.Range("E2").Value = '  -3.20%  '
$ws.Range("D3").Value = '1.855.98'
$ws.Range("E3").Value = '  -4.04%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").Value = '323.88'
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '0.4527'
$ws.Range("E7").Value = '  -4.63%  '
$ws.Range("D8").Value = '0.3866'
$ws.Range("E8").Value = '  -4.94%  '
$ws.Range("D9").Value = '48.54'
$ws.Range("E9").Value = '  -9.11%  '
$ws.Range("D10").Value = '0.07923'
$ws.Range("E10").Value = '  -6.98%  '
$ws.Range("D11").Value = '1.015'
$ws.Range("E11").Value = '  -3.61%  '
$ws.Range("D12").Value = '21.36'
$ws.Range("E12").Value = '  -4.36%  '
$ws.Range("D13").Value = '1.872.45'
$ws.Range("E13").Value = '  -5.05%  '
$ws.Range("D14").Value = '5.913'
$ws.Range("E14").Value = '  -3.71%  '
$ws.Range("D15").Value = '7.116'
$ws.Range("E15").Value = '  -5.74%  '
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '85.86'
$ws.Range("E17").Value = '  -4.74%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = '0.00001030'
$ws.Range("E18").Value = '  -3.74%  '
$ws.Range("D19").Value = '0.06522'
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  -6.38%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = '5.544'
$ws.Range("E22").Value = '  -4.23%  '
$ws.Range("D23").Value = '27.339.69'
$ws.Range("E23").Value = '  -3.27%  '
$ws.Range("D24").Value = '10.84'
$ws.Range("E24").Value = '  -5.44%  '
$ws.Range("D25").Value = '2.279'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = '2.085.65'
$ws.Range("E26").Value = '  -5.26%  '
$ws.Range("D27").Value = '153.85'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").Value = '19.92'
$ws.Range("E28").Value = '  -1.36%  '
$ws.Range("D29").Value = '2.074'
$ws.Range("E29").Value = '  -4.75%  '
$ws.Range("D30").Value = '5.431'
$ws.Range("E30").Value = '  -6.18%  '
$ws.Range("D31").Value = '120.98'
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("D32").Value = '1.481'
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("D33").Value = '0.09289'
$ws.Range("E33").Value = '  -3.41%  '
$ws.Range("D34").Value = '0.9349'
$ws.Range("E34").Value = '  -5.05%  '
$ws.Range("D35").Value = '3.599'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").Value = '5.264'
$ws.Range("E36").Value = '  -6.04%  '
$ws.Range("D37").Value = '1.229'
$ws.Range("E37").Value = '  -1.16%  '
$ws.Range("D38").Value = '0.02226'
$ws.Range("E38").Value = '  -4.47%  '
$ws.Range("D39").Value = '0.05990'
$ws.Range("E39").Value = '  -3.10%  '
$ws.Range("D40").Value = '8.180'
$ws.Range("E40").Value = '  -12.08%  '
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = '0.5910'
$ws.Range("E42").Value = '  -4.72%  '
$ws.Range("D43").Value = '0.1906'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '10.11'
$ws.Range("E44").Value = '  -9.48%  '
$ws.Range("D45").Value = '1.276'
$ws.Range("E45").Value = '  -3.17%  '
$ws.Range("D46").Value = '0.5623'
$ws.Range("E46").Value = '  -5.09%  '
$ws.Range("D47").Value = '12.06'
$ws.Range("E47").Value = '  -6.10%  '
$ws.Range("D48").Value = '3.376'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '1.920'
$ws.Range("E49").Value = '  -6.31%  '
$ws.Range("D50").Value = '0.06766'
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = '108.48'
$ws.Range("E51").Value = '  -1.46%  '
